$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.857.98'
$ws.Range('E2').Value = '  +3.00%  '
$ws.Range('D3').Value = '2.092.55'
$ws.Range('E3').Value = '  +2.43%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.72'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.616'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.75%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.51'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.71%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +2.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0837'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.56%  '
$ws.Range('E11').Value = '  +0.21%  '
$ws.Range('D12').Value = '2.402.41'
$ws.Range('E12').Value = '  +2.41%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.99'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.91'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.795'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.21%  '
$ws.Range('E16').Value = '  -0.44%  '
$ws.Range('D17').Value = '2.094.18'
$ws.Range('E17').Value = '  +2.51%  '
$ws.Range('D18').Value = '38.702.89'
$ws.Range('E18').Value = '  +2.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '71.64'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.36%  '
$ws.Range('E20').Value = '  +2.30%  '
$ws.Range('D21').Value = '0.0₃0836'
$ws.Range('E21').Value = '  +1.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '227.39'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.43%  '
$ws.Range('E23').Value = '  -0.45%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.38'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('E25').Value = '  +3.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '171.28'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.50'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.15%  '
$ws.Range('E28').Value = '  +9.52%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.48'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +15.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.18'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.31%  '
$ws.Range('E31').Value = '  +1.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.35'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.39%  '
$ws.Range('E33').Value = '  +3.10%  '
$ws.Range('E34').Value = '  +4.46%  '
$ws.Range('E35').Value = '  +1.75%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.45'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.39'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.03%  '
$ws.Range('E38').Value = '  +3.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.21%  '
$ws.Range('E40').Value = '  +0.53%  '
$ws.Range('D41').Value = '1.540.67'
$ws.Range('E41').Value = '  +0.60%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '100.85'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.35%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0225'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.45%  '
$ws.Range('E44').Value = '  -0.76%  '
$ws.Range('E45').Value = '  +3.76%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.64'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +8.87%  '
$ws.Range('E47').Value = '  +1.34%  '
$ws.Range('E48').Value = '  -1.12%  '
$ws.Range('E49').Value = '  +3.05%  '
$ws.Range('E50').Value = '  +1.17%  '
$ws.Range('D51').Value = '2.288.41'
$ws.Range('E51').Value = '  +2.41%  '
